# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the Pina dataset block
# (rows 114-115), pushing the existing rows 114-144 down to 146.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 114 (shifts old 114..144 down to 116..146)
$ws.Rows.Item(114).Insert()
$ws.Rows.Item(114).Insert()

# New row 114: "Primera" quality
$ws.Range("A114").Value = 7
$ws.Range("B114").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C114").Value = "Ñuble"
$ws.Range("D114").Value = 44468
$ws.Range("E114").Value = 16
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100108
$ws.Range("H114").Value = "Tropicales y subtropicales"
$ws.Range("I114").Value = 100108005
$ws.Range("J114").Value = "Piña"
$ws.Range("K114").Value = "Caramelo"
$ws.Range("L114").Value = "Primera"
$ws.Range("M114").Value = 60
$ws.Range("N114").Value = 20000
$ws.Range("O114").Value = 21000
$ws.Range("P114").Value = 20500
$ws.Range("Q114").Value = '$/caja 12 unidades'
$ws.Range("R114").Value = "Ecuador"
$ws.Range("S114").Value = 1708
$ws.Range("T114").Value = 12

# New row 115: "Segunda" quality
$ws.Range("A115").Value = 7
$ws.Range("B115").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C115").Value = "Ñuble"
$ws.Range("D115").Value = 44468
$ws.Range("E115").Value = 16
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100108
$ws.Range("H115").Value = "Tropicales y subtropicales"
$ws.Range("I115").Value = 100108005
$ws.Range("J115").Value = "Piña"
$ws.Range("K115").Value = "Caramelo"
$ws.Range("L115").Value = "Segunda"
$ws.Range("M115").Value = 60
$ws.Range("N115").Value = 20000
$ws.Range("O115").Value = 21000
$ws.Range("P115").Value = 20500
$ws.Range("Q115").Value = '$/caja 14 unidades'
$ws.Range("R115").Value = "Ecuador"
$ws.Range("S115").Value = 1464
$ws.Range("T115").Value = 14
